$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.497.03'
$ws.Range('E2').Value = '  -4.38%  '
$ws.Range('D3').Value = '2.928.99'
$ws.Range('E3').Value = '  -2.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '548.20'
$ws.Range('E5').Value = '  -4.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.31'
$ws.Range('E6').Value = '  +4.24%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.60%  '
$ws.Range('D9').Value = '2.923.35'
$ws.Range('E9').Value = '  -2.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.127'
$ws.Range('E10').Value = '  -3.51%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '4.77'
$ws.Range('E11').Value = '  -5.79%  '
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.87'
$ws.Range('E14').Value = '  +1.27%  '
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('D16').Value = '3.410.92'
$ws.Range('E16').Value = '  -2.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.85'
$ws.Range('E17').Value = '  +6.69%  '
$ws.Range('D18').Value = '2.926.08'
$ws.Range('E18').Value = '  -2.14%  '
$ws.Range('D19').Value = '57.536.60'
$ws.Range('E19').Value = '  -4.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '417.35'
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.18'
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.686'
$ws.Range('E22').Value = '  +2.80%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.97'
$ws.Range('E23').Value = '  -1.00%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.99'
$ws.Range('E24').Value = '  +0.83%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '79.83'
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.47'
$ws.Range('E28').Value = '  -2.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.46'
$ws.Range('E29').Value = '  +3.14%  '
$ws.Range('E30').Value = '  +2.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.20'
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.96'
$ws.Range('E32').Value = '  -2.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0970'
$ws.Range('E33').Value = '  +2.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.65'
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.940'
$ws.Range('E35').Value = '  +0.66%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.07'
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('B37').Value = 'Cosmos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.74'
$ws.Range('E37').Value = '  +4.19%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '47.93'
$ws.Range('E38').Value = '  -4.37%  '
$ws.Range('D39').Value = '0.0₃0679'
$ws.Range('E39').Value = '  +2.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.54'
$ws.Range('E40').Value = '  +3.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.108'
$ws.Range('E41').Value = '  -0.59%  '
$ws.Range('E42').Value = '  -2.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '374.74'
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '2.679.12'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.26'
$ws.Range('E47').Value = '  +1.63%  '
$ws.Range('E48').Value = '  +1.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.97'
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.14'
$ws.Range('E50').Value = '  -1.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.00'
$ws.Range('E51').Value = '  -0.32%  '
